# Improve cell segmentation algorithm: use histogram for identifying
# cells with no inclusions. This updates the Number_of_Inclusions (col B)
# and recomputed Number_of_Inclusions_per_Nucleus (col D) values for the
# rows affected by the refined segmentation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new Number_of_Inclusions (column B) value
$updates = @{
    2  = 20
    3  = 4
    4  = 6
    5  = 2
    6  = 41
    7  = 5
    8  = 10
    12 = 89
    14 = 2
    19 = 2
    22 = 2
    23 = 8
    27 = 1
    28 = 0
    29 = 2
}

foreach ($row in $updates.Keys) {
    $newInclusions = $updates[$row]

    # Column B: Number_of_Inclusions
    $ws.Cells.Item($row, 2).Value = $newInclusions

    # Column C: Number_of_Nuclei (unchanged, read for recompute)
    $nuclei = $ws.Cells.Item($row, 3).Value()

    # Column D: Number_of_Inclusions_per_Nucleus = Inclusions / Nuclei
    # Excel stores doubles at ~15-16 significant-digit precision, so the
    # ratio is round-tripped through a 16-significant-digit string
    # (matches Excel's own floating-point display/storage behavior)
    # before being written back to the cell.
    if ($nuclei -ne 0) {
        $ratio = $newInclusions / $nuclei
        $ws.Cells.Item($row, 4).Value = $ratio.ToString("G16")
    } else {
        $ws.Cells.Item($row, 4).Value = 0
    }
}
